$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("C1").Value = "White"
$ws.Range("E1").Value = "Black or African"
$ws.Range("G1").Value = "American Indian & Alaska"
$ws.Range("H1").Value = "Asian Count"
$ws.Range("I1").Value = "Asian"
$ws.Range("J1").Value = "Native Hawaiian Count"
$ws.Range("K1").Value = "Native Hawaiian"
$ws.Range("M1").Value = "Hispanic or Latino"
$ws.Range("O1").Value = "Others"
$ws.Range("Q1").Value = "Two or More Ethnicity"

$ws.Range("N5").Select()
